# Adds a new user row ("hoal") to the users sheet, mirroring the existing
# rows (Username/password/email/name/surname1/surname2/phone).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A6").Value = "hoal"
$ws.Range("B6").Value = "hoal"
$ws.Range("C6").Value = "hoal@gmail.com"
$ws.Range("D6").Value = "hoal"
$ws.Range("E6").Value = "hoal"
$ws.Range("F6").Value = "hoal"
$ws.Range("G6").Value = 555999111
